$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text without Excel coercing it to a number/percent,
# and without leaving a residual custom style behind on the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '305.47'
Set-TextValue $ws.Range('E2') '0.65%'
Set-TextValue $ws.Range('D3') '37.17'
Set-TextValue $ws.Range('E3') '6.26%'
Set-TextValue $ws.Range('D4') '5.014'
Set-TextValue $ws.Range('E4') '-3.01%'
Set-TextValue $ws.Range('D5') '0.07900'
Set-TextValue $ws.Range('E5') '0.90%'
Set-TextValue $ws.Range('D6') '2.211'
Set-TextValue $ws.Range('E6') '-5.39%'
Set-TextValue $ws.Range('D7') '8.024'
Set-TextValue $ws.Range('E7') '-0.32%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.9243'
Set-TextValue $ws.Range('E8') '-0.06%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D9') '0.09708'
Set-TextValue $ws.Range('E9') '-3.85%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D10') '0.1885'
Set-TextValue $ws.Range('E10') '3.25%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D11') '0.08621'
Set-TextValue $ws.Range('E11') '0.87%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D12') '0.03696'
Set-TextValue $ws.Range('E12') '7.99%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D13') '0.09979'
Set-TextValue $ws.Range('E13') '0.59%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D14') '0.001469'
Set-TextValue $ws.Range('E14') '-0.90%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D15') '0.005649'
Set-TextValue $ws.Range('E15') '-1.00%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D16') '3.467'
Set-TextValue $ws.Range('E16') '-0.14%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D17') '4.023'
Set-TextValue $ws.Range('E17') '1.06%'
Set-TextValue $ws.Range('E18') '6.65%'
Set-TextValue $ws.Range('D19') '0.3414'
Set-TextValue $ws.Range('E19') '-0.37%'
Set-TextValue $ws.Range('E20') '-0.68%'
Set-TextValue $ws.Range('E21') '4.39%'
Set-TextValue $ws.Range('D22') '0.2201'
Set-TextValue $ws.Range('E22') '-1.08%'
Set-TextValue $ws.Range('D23') '0.04563'
Set-TextValue $ws.Range('E23') '-1.91%'
Set-TextValue $ws.Range('E24') '1.11%'
Set-TextValue $ws.Range('D25') '0.004474'
Set-TextValue $ws.Range('E25') '3.23%'
Set-TextValue $ws.Range('D26') '0.0001397'
Set-TextValue $ws.Range('E26') '7.43%'
Set-TextValue $ws.Range('E27') '39.61%'
Set-TextValue $ws.Range('D39') '0.01844'
Set-TextValue $ws.Range('E39') '5.28%'
Set-TextValue $ws.Range('D40') '0.04802'
Set-TextValue $ws.Range('E40') '1.16%'
Set-TextValue $ws.Range('D41') '0.008125'
Set-TextValue $ws.Range('E41') '4.66%'
Set-TextValue $ws.Range('E42') '-1.13%'
Set-TextValue $ws.Range('D43') '0.007564'
Set-TextValue $ws.Range('E43') '-14.65%'
Set-TextValue $ws.Range('D44') '0.002205'
Set-TextValue $ws.Range('E44') '-0.25%'
Set-TextValue $ws.Range('D45') '0.01006'
Set-TextValue $ws.Range('E45') '0.86%'
Set-TextValue $ws.Range('D46') '0.00006258'
Set-TextValue $ws.Range('E46') '2.96%'
Set-TextValue $ws.Range('E47') '-0.09%'
Set-TextValue $ws.Range('E48') '0.02%'
Set-TextValue $ws.Range('E49') '390.30%'
Set-TextValue $ws.Range('D50') '0.001722'
Set-TextValue $ws.Range('E50') '-36.08%'
Set-TextValue $ws.Range('D51') '0.00002101'
Set-TextValue $ws.Range('E51') '-0.09%'
